$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates are OLE Automation day-count serials, matching
# the existing Date column's format) appended below the current last
# row (row 48).
$dates = @(46031, 46036, 46034, 46035)
$counts = @(54, 91, 75, 102)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 49 + $i

    # Copy the date-format (s="1" / numFmtId 14) from the cell directly
    # above so the new rows keep the same date styling instead of Excel
    # minting a fresh custom number format.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Range("A" + $row).Value = $dates[$i]
    $ws.Range("B" + $row).Value = $counts[$i]
}

$ws.Range("A50:B52").Select() | Out-Null
